{"js": "// Replace the math-fact answers in the single results table.\n// The mapping below is keyed by the *old* cell text (each is unique in the\n// document) -> the *new* cell text, taken directly from the authoritative\n// diff. We read the table's full `values` grid, remap every cell that has\n// an entry in the mapping, and write the grid back in one shot so existing\n// run formatting (font/size) on each cell is preserved.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document, found none.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst mapping = {\n  \"0+75=75\": \"89-24=65\",\n  \"75-53=22\": \"18+7=25\",\n  \"73-48=25\": \"19+62=81\",\n  \"76-14=62\": \"95-34=61\",\n  \"1+4=5\": \"57-29=28\",\n  \"47-22=25\": \"41-6=35\",\n  \"20+18=38\": \"40+34=74\",\n  \"22+44=66\": \"27+31=58\",\n  \"64-52=12\": \"52-7=45\",\n  \"96-4=92\": \"46-17=29\",\n  \"21+0=21\": \"86+0=86\",\n  \"64-48=16\": \"54+37=91\",\n  \"25-17=8\": \"95-16=79\",\n  \"43+12=55\": \"69-58=11\",\n  \"8-7=1\": \"90-71=19\",\n  \"74-43=31\": \"23-18=5\",\n  \"68-24=44\": \"73-64=9\",\n  \"23-0=23\": \"78-51=27\",\n  \"10+19=29\": \"46-13=33\",\n  \"33-8=25\": \"11+10=21\",\n  \"24+4=28\": \"59+36=95\",\n  \"26+32=58\": \"45-32=13\",\n  \"22+20=42\": \"70-35=35\",\n  \"91-73=18\": \"65-21=44\",\n  \"52-15=37\": \"22-1=21\",\n  \"91-66=25\": \"44+27=71\",\n  \"40-35=5\": \"86-51=35\",\n  \"91-36=55\": \"8+0=8\",\n  \"64-23=41\": \"82-77=5\",\n  \"59+19=78\": \"68-50=18\",\n  \"27+45=72\": \"0+84=84\",\n  \"49+11=60\": \"99-92=7\",\n  \"90-51=39\": \"36-23=13\",\n  \"9+68=77\": \"78-51=27\",\n  \"65-27=38\": \"30+10=40\",\n  \"71-4=67\": \"40-1=39\",\n  \"1+73=74\": \"58+13=71\",\n  \"14+13=27\": \"40+26=66\",\n  \"62-41=21\": \"52-40=12\",\n  \"92-25=67\": \"44-18=26\",\n  \"92+2=94\": \"79+3=82\",\n  \"82-71=11\": \"93-56=37\",\n  \"26+72=98\": \"96-31=65\",\n  \"41-27=14\": \"60-19=41\",\n  \"86+5=91\": \"17+62=79\",\n  \"19+57=76\": \"41-21=20\",\n  \"40+31=71\": \"77-9=68\",\n  \"41-16=25\": \"20+70=90\",\n  \"20+53=73\": \"33+32=65\",\n  \"70-23=47\": \"52-28=24\",\n  \"96-0=96\": \"78-4=74\",\n  \"42+40=82\": \"99-31=68\",\n  \"86-23=63\": \"59-42=17\",\n  \"95-64=31\": \"94-92=2\",\n  \"14+11=25\": \"62-1=61\",\n  \"18+73=91\": \"28+8=36\",\n  \"59+37=96\": \"92-36=56\",\n  \"93-89=4\": \"28+70=98\",\n  \"86-83=3\": \"85-82=3\",\n  \"37+6=43\": \"88-47=41\",\n  \"74+23=97\": \"58+6=64\",\n  \"16+52=68\": \"4+51=55\",\n  \"21+52=73\": \"97-17=80\",\n  \"3+30=33\": \"3+88=91\",\n  \"92-53=39\": \"48+9=57\",\n  \"5+91=96\": \"78-63=15\",\n  \"78-11=67\": \"15+18=33\",\n  \"1+22=23\": \"78-70=8\",\n  \"63+20=83\": \"25-21=4\",\n  \"25+16=41\": \"37+59=96\",\n  \"80-75=5\": \"27-15=12\",\n  \"60-33=27\": \"41-31=10\",\n  \"28+27=55\": \"86-11=75\",\n  \"65+34=99\": \"19+50=69\",\n  \"85+5=90\": \"86-7=79\",\n  \"26+59=85\": \"4+18=22\",\n  \"73-67=6\": \"26-6=20\",\n  \"33-23=10\": \"17+70=87\",\n  \"63+25=88\": \"83-12=71\",\n  \"57-22=35\": \"30+59=89\",\n  \"19+43=62\": \"85+2=87\",\n  \"66-31=35\": \"55-28=27\",\n  \"2+21=23\": \"78-73=5\",\n  \"57-38=19\": \"36-32=4\",\n  \"18+26=44\": \"27+67=94\",\n  \"58+22=80\": \"34+49=83\",\n  \"59+39=98\": \"75+7=82\",\n  \"89-50=39\": \"7-6=1\",\n  \"63-39=24\": \"47-15=32\",\n  \"88+2=90\": \"20+59=79\",\n  \"49-42=7\": \"68-68=0\",\n  \"50+24=74\": \"77+7=84\",\n  \"60-45=15\": \"17+78=95\",\n  \"59+27=86\": \"99-56=43\",\n  \"79-43=36\": \"2+18=20\",\n  \"98-69=29\": \"38-15=23\",\n  \"66-22=44\": \"84-24=60\",\n  \"81+15=96\": \"44+8=52\",\n  \"32+66=98\": \"82-5=77\",\n  \"33+17=50\": \"56+39=95\",\n};\n\nconst newValues = table.values.map(row =>\n  row.map(cellText => (Object.prototype.hasOwnProperty.call(mapping, cellText) ? mapping[cellText] : cellText))\n);\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Update the answers in the single \"addition/subtraction within 100\" results\n# table. Every old equation string is unique in the document, so a series of\n# exact, whole-word Find/Replace passes over the document body safely retargets\n# exactly one table cell each, leaving all paragraph/run formatting untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('0+75=75', '89-24=65'),\n    @('75-53=22', '18+7=25'),\n    @('73-48=25', '19+62=81'),\n    @('76-14=62', '95-34=61'),\n    @('1+4=5', '57-29=28'),\n    @('47-22=25', '41-6=35'),\n    @('20+18=38', '40+34=74'),\n    @('22+44=66', '27+31=58'),\n    @('64-52=12', '52-7=45'),\n    @('96-4=92', '46-17=29'),\n    @('21+0=21', '86+0=86'),\n    @('64-48=16', '54+37=91'),\n    @('25-17=8', '95-16=79'),\n    @('43+12=55', '69-58=11'),\n    @('8-7=1', '90-71=19'),\n    @('74-43=31', '23-18=5'),\n    @('68-24=44', '73-64=9'),\n    @('23-0=23', '78-51=27'),\n    @('10+19=29', '46-13=33'),\n    @('33-8=25', '11+10=21'),\n    @('24+4=28', '59+36=95'),\n    @('26+32=58', '45-32=13'),\n    @('22+20=42', '70-35=35'),\n    @('91-73=18', '65-21=44'),\n    @('52-15=37', '22-1=21'),\n    @('91-66=25', '44+27=71'),\n    @('40-35=5', '86-51=35'),\n    @('91-36=55', '8+0=8'),\n    @('64-23=41', '82-77=5'),\n    @('59+19=78', '68-50=18'),\n    @('27+45=72', '0+84=84'),\n    @('49+11=60', '99-92=7'),\n    @('90-51=39', '36-23=13'),\n    @('9+68=77', '78-51=27'),\n    @('65-27=38', '30+10=40'),\n    @('71-4=67', '40-1=39'),\n    @('1+73=74', '58+13=71'),\n    @('14+13=27', '40+26=66'),\n    @('62-41=21', '52-40=12'),\n    @('92-25=67', '44-18=26'),\n    @('92+2=94', '79+3=82'),\n    @('82-71=11', '93-56=37'),\n    @('26+72=98', '96-31=65'),\n    @('41-27=14', '60-19=41'),\n    @('86+5=91', '17+62=79'),\n    @('19+57=76', '41-21=20'),\n    @('40+31=71', '77-9=68'),\n    @('41-16=25', '20+70=90'),\n    @('20+53=73', '33+32=65'),\n    @('70-23=47', '52-28=24'),\n    @('96-0=96', '78-4=74'),\n    @('42+40=82', '99-31=68'),\n    @('86-23=63', '59-42=17'),\n    @('95-64=31', '94-92=2'),\n    @('14+11=25', '62-1=61'),\n    @('18+73=91', '28+8=36'),\n    @('59+37=96', '92-36=56'),\n    @('93-89=4', '28+70=98'),\n    @('86-83=3', '85-82=3'),\n    @('37+6=43', '88-47=41'),\n    @('74+23=97', '58+6=64'),\n    @('16+52=68', '4+51=55'),\n    @('21+52=73', '97-17=80'),\n    @('3+30=33', '3+88=91'),\n    @('92-53=39', '48+9=57'),\n    @('5+91=96', '78-63=15'),\n    @('78-11=67', '15+18=33'),\n    @('1+22=23', '78-70=8'),\n    @('63+20=83', '25-21=4'),\n    @('25+16=41', '37+59=96'),\n    @('80-75=5', '27-15=12'),\n    @('60-33=27', '41-31=10'),\n    @('28+27=55', '86-11=75'),\n    @('65+34=99', '19+50=69'),\n    @('85+5=90', '86-7=79'),\n    @('26+59=85', '4+18=22'),\n    @('73-67=6', '26-6=20'),\n    @('33-23=10', '17+70=87'),\n    @('63+25=88', '83-12=71'),\n    @('57-22=35', '30+59=89'),\n    @('19+43=62', '85+2=87'),\n    @('66-31=35', '55-28=27'),\n    @('2+21=23', '78-73=5'),\n    @('57-38=19', '36-32=4'),\n    @('18+26=44', '27+67=94'),\n    @('58+22=80', '34+49=83'),\n    @('59+39=98', '75+7=82'),\n    @('89-50=39', '7-6=1'),\n    @('63-39=24', '47-15=32'),\n    @('88+2=90', '20+59=79'),\n    @('49-42=7', '68-68=0'),\n    @('50+24=74', '77+7=84'),\n    @('60-45=15', '17+78=95'),\n    @('59+27=86', '99-56=43'),\n    @('79-43=36', '2+18=20'),\n    @('98-69=29', '38-15=23'),\n    @('66-22=44', '84-24=60'),\n    @('81+15=96', '44+8=52'),\n    @('32+66=98', '82-5=77'),\n    @('33+17=50', '56+39=95')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    # wdFindContinue = 1, wdReplaceOne = 1 (only the single exact match per pair)\n    $find.Execute(\n        $oldText,   # FindText\n        $false,     # MatchCase\n        $true,      # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        1           # Replace (wdReplaceOne)\n    ) | Out-Null\n}\n"}
